# Apply the "Updated cryptos list" refresh (Thu Nov 21 14:51:33 UTC 2024 GitHub Actions run).
# For each changed row: Price (column D) and Volume/1h change (column E) are refreshed with
# newly scraped figures; rows 30/31 (Cronos vs. Stellar) additionally swap position, so their
# Coin name / Link / Price / Volume cells are updated in full.
#
# Several Price values (column D) look like plain numbers (e.g. "0.999", "244.90") but must
# stay as literal TEXT, matching the sheet's existing inline-string cells (e.g. trailing zeros
# such as "244.90" must be preserved exactly, not normalised to 244.9). A leading apostrophe
# forces Excel to store these as text (quote-prefixed) instead of re-parsing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '97.457.38'
$ws.Range("E2").Value = '  +3.03%  '

$ws.Range("D3").Value = '3.353.81'
$ws.Range("E3").Value = '  +7.62%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").Value = "'244.90"
$ws.Range("E5").Value = '  +2.18%  '

$ws.Range("D6").Value = "'625.63"
$ws.Range("E6").Value = '  +1.58%  '

$ws.Range("D7").Value = "'1.13"
$ws.Range("E7").Value = '  -0.30%  '

$ws.Range("D8").Value = "'0.389"
$ws.Range("E8").Value = '  -0.81%  '

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("D10").Value = '3.359.28'
$ws.Range("E10").Value = '  +7.83%  '

$ws.Range("D11").Value = "'0.792"
$ws.Range("E11").Value = '  -4.32%  '

$ws.Range("E12").Value = '  +0.63%  '

$ws.Range("D13").Value = '96.985.87'
$ws.Range("E13").Value = '  +2.95%  '

$ws.Range("D14").Value = "'0.0000248"
$ws.Range("E14").Value = '  +0.97%  '

$ws.Range("D15").Value = "'35.27"
$ws.Range("E15").Value = '  +1.45%  '

$ws.Range("D16").Value = '3.963.69'

$ws.Range("D17").Value = "'5.53"
$ws.Range("E17").Value = '  +2.14%  '

$ws.Range("D18").Value = '3.338.83'
$ws.Range("E18").Value = '  +6.80%  '

$ws.Range("D19").Value = "'3.58"
$ws.Range("E19").Value = '  -3.76%  '

$ws.Range("D20").Value = "'15.22"
$ws.Range("E20").Value = '  +1.50%  '

$ws.Range("D21").Value = "'491.19"
$ws.Range("E21").Value = '  +8.79%  '

$ws.Range("D22").Value = "'0.0000210"
$ws.Range("E22").Value = '  +4.74%  '

$ws.Range("D23").Value = "'5.87"
$ws.Range("E23").Value = '  -1.45%  '

$ws.Range("D24").Value = "'9.33"
$ws.Range("E24").Value = '  +3.44%  '

$ws.Range("D25").Value = "'5.71"
$ws.Range("E25").Value = '  +0.90%  '

$ws.Range("D26").Value = "'88.45"
$ws.Range("E26").Value = '  +2.22%  '

$ws.Range("D27").Value = "'12.17"
$ws.Range("E27").Value = '  +0.53%  '

$ws.Range("D28").Value = '3.537.59'
$ws.Range("E28").Value = '  +7.46%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = "'0.243"
$ws.Range("E30").Value = '  -7.17%  '

$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").Value = "'0.181"
$ws.Range("E31").Value = '  +0.38%  '

$ws.Range("E32").Value = '  -0.11%  '

$ws.Range("E33").Value = '  -4.32%  '

$ws.Range("D34").Value = "'9.36"
$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").Value = "'27.72"
$ws.Range("E35").Value = '  +5.66%  '

$ws.Range("D36").Value = "'7.45"
$ws.Range("E36").Value = '  -5.79%  '

$ws.Range("E37").Value = '  -6.35%  '

$ws.Range("D38").Value = "'1.94"
$ws.Range("E38").Value = '  +1.14%  '

$ws.Range("D39").Value = "'499.58"
$ws.Range("E39").Value = '  +4.21%  '

$ws.Range("D40").Value = "'24.64"
$ws.Range("E40").Value = '  +2.83%  '

$ws.Range("D41").Value = "'0.452"
$ws.Range("E41").Value = '  -0.76%  '

$ws.Range("D42").Value = "'1.28"
$ws.Range("E42").Value = '  -0.50%  '

$ws.Range("D43").Value = "'0.811"
$ws.Range("E43").Value = '  +16.70%  '

$ws.Range("D44").Value = "'3.29"
$ws.Range("E44").Value = '  -0.38%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").Value = "'3.45"
$ws.Range("E46").Value = '  -7.86%  '

$ws.Range("D47").Value = "'159.74"
$ws.Range("E47").Value = '  -0.65%  '

$ws.Range("D48").Value = "'1.94"
$ws.Range("E48").Value = '  +3.68%  '

$ws.Range("D49").Value = "'4.57"
$ws.Range("E49").Value = '  +2.85%  '

$ws.Range("D50").Value = "'45.29"
$ws.Range("E50").Value = '  +3.24%  '

$ws.Range("E51").Value = '  +2.18%  '
